$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: value update only
$ws.Range("B2").Value = 0.1333333333333333

# Row 4: value update only
$ws.Range("B4").Value = 0

# Row 5: value update only
$ws.Range("B5").Value = 0

# Rows 10-14: labels shift up by one (wrap-around) and values updated
$ws.Range("A10").Value = "Disparity Generators"
$ws.Range("B10").Value = 0

$ws.Range("A11").Value = "Disparity Load"
$ws.Range("B11").Value = 0.5288308596787999

$ws.Range("A12").Value = "Disparity Trafo"
$ws.Range("B12").Value = 0

$ws.Range("A13").Value = "Disparity Lines"
$ws.Range("B13").Value = 0.7993489242141023

$ws.Range("A14").Value = "Overall 70% Redundancy"
$ws.Range("B14").Value = 0.9052132701421801
